$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibition) - column F "想去人数"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 32
$ws1.Cells.Item(3, 6).Value = 8826
$ws1.Cells.Item(4, 6).Value = 1950
$ws1.Cells.Item(5, 6).Value = 6563
$ws1.Cells.Item(6, 6).Value = 166
$ws1.Cells.Item(8, 6).Value = 588
$ws1.Cells.Item(11, 6).Value = 65
$ws1.Cells.Item(15, 6).Value = 15
$ws1.Cells.Item(16, 6).Value = 8720
$ws1.Cells.Item(20, 6).Value = 115
$ws1.Cells.Item(21, 6).Value = 1825
$ws1.Cells.Item(26, 6).Value = 65
$ws1.Cells.Item(28, 6).Value = 196
$ws1.Cells.Item(29, 6).Value = 1020
$ws1.Cells.Item(30, 6).Value = 22
$ws1.Cells.Item(31, 6).Value = 38
$ws1.Cells.Item(32, 6).Value = 26
$ws1.Cells.Item(33, 6).Value = 25
$ws1.Cells.Item(34, 6).Value = 2192
$ws1.Cells.Item(35, 6).Value = 866
$ws1.Cells.Item(36, 6).Value = 511
$ws1.Cells.Item(40, 6).Value = 236
$ws1.Cells.Item(41, 6).Value = 170
$ws1.Cells.Item(43, 6).Value = 350
$ws1.Cells.Item(46, 6).Value = 65
$ws1.Cells.Item(47, 6).Value = 3984

# Sheet: 演出 (Performance) - column F
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 402
$ws2.Cells.Item(24, 6).Value = 67

# Sheet: 本地生活 (Local Life) - column F
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(3, 6).Value = 715
$ws3.Cells.Item(4, 6).Value = 321

# Sheet: 全部类型 (All Types) - column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(3, 6).Value = 715
$ws4.Cells.Item(4, 6).Value = 402
$ws4.Cells.Item(5, 6).Value = 8826
$ws4.Cells.Item(7, 6).Value = 321
$ws4.Cells.Item(8, 6).Value = 1950
$ws4.Cells.Item(9, 6).Value = 6563
$ws4.Cells.Item(10, 6).Value = 166
$ws4.Cells.Item(13, 6).Value = 588
$ws4.Cells.Item(17, 6).Value = 65
$ws4.Cells.Item(19, 6).Value = 15
$ws4.Cells.Item(20, 6).Value = 8720
$ws4.Cells.Item(23, 6).Value = 115
$ws4.Cells.Item(24, 6).Value = 1825
$ws4.Cells.Item(28, 6).Value = 65
$ws4.Cells.Item(30, 6).Value = 196
$ws4.Cells.Item(31, 6).Value = 1020
$ws4.Cells.Item(32, 6).Value = 38
$ws4.Cells.Item(33, 6).Value = 25
$ws4.Cells.Item(34, 6).Value = 2192
$ws4.Cells.Item(35, 6).Value = 866
$ws4.Cells.Item(36, 6).Value = 511
$ws4.Cells.Item(40, 6).Value = 236
$ws4.Cells.Item(42, 6).Value = 170
$ws4.Cells.Item(45, 6).Value = 3984
